# Edit script generated to apply the "data up to 7th" commit
# Adds column E (American Samoa) values for rows 117-122,
# fills full data rows for 211-216, and appends new date label rows 217-221.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill previously-missing column E (American Samoa) values for rows 117-122 ---
$ws.Range("E117").Value = 0.11989087301587
$ws.Range("E118").Value = 0.1210729746444
$ws.Range("E119").Value = 0.11715888278388
$ws.Range("E120").Value = 0.1140380713689
$ws.Range("E121").Value = 0.11938431938432
$ws.Range("E122").Value = 0.12386363636364

# --- Fill data rows 211-216 (columns B:BE, column E left blank) ---
$row211 = New-Object "object[,]" 1,56
$row211[0,0] = 0.051696994500759
$row211[0,1] = 0.07625683015227699
$row211[0,2] = 0.07085532900848
$row211[0,3] = $null
$row211[0,4] = 0.04404557266281
$row211[0,5] = 0.049105169044618
$row211[0,6] = 0.051758253159025
$row211[0,7] = 0.046945989094865
$row211[0,8] = 0.043905507439192
$row211[0,9] = 0.042771435587403
$row211[0,10] = 0.050419432468566
$row211[0,11] = 0.054364373980746
$row211[0,12] = 0.058680272055385
$row211[0,13] = 0.023611326051445
$row211[0,14] = 0.07585935513485099
$row211[0,15] = 0.048132590130703
$row211[0,16] = 0.053085858786677
$row211[0,17] = 0.058726580659376
$row211[0,18] = 0.062408184442694
$row211[0,19] = 0.054507255684291
$row211[0,20] = 0.055057282939854
$row211[0,21] = 0.048224630508421
$row211[0,22] = 0.043306993504741
$row211[0,23] = 0.051985070815827
$row211[0,24] = 0.04740499263408
$row211[0,25] = 0.044763809796584
$row211[0,26] = 0.072209816566874
$row211[0,27] = 0.06612266070954199
$row211[0,28] = 0.07189134440184899
$row211[0,29] = 0.057366447723123
$row211[0,30] = 0.052388695794005
$row211[0,31] = 0.068057500374391
$row211[0,32] = 0.082812683650829
$row211[0,33] = 0.048752536722536
$row211[0,34] = 0.041786915853282
$row211[0,35] = 0.044876666638354
$row211[0,36] = 0.041090261763591
$row211[0,37] = 0.044754490371372
$row211[0,38] = 0.056952059376988
$row211[0,39] = 0.067974558290786
$row211[0,40] = 0.042280091903767
$row211[0,41] = 0.051659138794791
$row211[0,42] = 0.07793119066847
$row211[0,43] = 0.048588828356334
$row211[0,44] = 0.053345396613048
$row211[0,45] = 0.06496672563651899
$row211[0,46] = 0.062520459748382
$row211[0,47] = 0.054917192317404
$row211[0,48] = 0.052617974668983
$row211[0,49] = 0.048300142643056
$row211[0,50] = 0.076457635330892
$row211[0,51] = 0.055280788001795
$row211[0,52] = 0.040970913842408
$row211[0,53] = 0.054420665774628
$row211[0,54] = 0.054883776013044
$row211[0,55] = 0.055005842462358
$ws.Range("B211:BE211").Value = $row211

$row212 = New-Object "object[,]" 1,56
$row212[0,0] = 0.059951993256956
$row212[0,1] = 0.040136079050407
$row212[0,2] = 0.023527786034028
$row212[0,3] = $null
$row212[0,4] = 0.035653063567919
$row212[0,5] = 0.034378335999486
$row212[0,6] = 0.037309606919259
$row212[0,7] = 0.032237255137754
$row212[0,8] = 0.046169573074369
$row212[0,9] = 0.030059465517806
$row212[0,10] = 0.02784903764383
$row212[0,11] = 0.025775016254533
$row212[0,12] = 0.06451213319409201
$row212[0,13] = 0.026019083222693
$row212[0,14] = 0.032493642517419
$row212[0,15] = 0.036548490555898
$row212[0,16] = 0.032594282851743
$row212[0,17] = 0.030704713533832
$row212[0,18] = 0.031254062449958
$row212[0,19] = 0.030667769818811
$row212[0,20] = 0.036149734980448
$row212[0,21] = 0.041223060795628
$row212[0,22] = 0.034661134396754
$row212[0,23] = 0.041212571069121
$row212[0,24] = 0.034515327943683
$row212[0,25] = 0.032205806775399
$row212[0,26] = 0.029999133744896
$row212[0,27] = 0.068363320548967
$row212[0,28] = 0.028572515099519
$row212[0,29] = 0.044216481802844
$row212[0,30] = 0.027399669561213
$row212[0,31] = 0.032787098225307
$row212[0,32] = 0.030963881037132
$row212[0,33] = 0.032294327061809
$row212[0,34] = 0.031819457037894
$row212[0,35] = 0.04339376571078
$row212[0,36] = 0.035101105394294
$row212[0,37] = 0.039320625811448
$row212[0,38] = 0.030327215563806
$row212[0,39] = 0.026911800976002
$row212[0,40] = 0.035421807917792
$row212[0,41] = 0.033495731766905
$row212[0,42] = 0.072342625453506
$row212[0,43] = 0.038492160385492
$row212[0,44] = 0.029753166084542
$row212[0,45] = 0.034018591926186
$row212[0,46] = 0.028182108951666
$row212[0,47] = 0.02884365978012
$row212[0,48] = 0.039534171006767
$row212[0,49] = 0.031712130636184
$row212[0,50] = 0.069777542337757
$row212[0,51] = 0.041950480641743
$row212[0,52] = 0.037244033230171
$row212[0,53] = 0.032487601257844
$row212[0,54] = 0.034520683785574
$row212[0,55] = 0.039935213308327
$ws.Range("B212:BE212").Value = $row212

$row213 = New-Object "object[,]" 1,56
$row213[0,0] = 0.054385154121895
$row213[0,1] = 0.032228053316519
$row213[0,2] = 0.02763487192516
$row213[0,3] = $null
$row213[0,4] = 0.03792468902438
$row213[0,5] = 0.037963869478382
$row213[0,6] = 0.038490797235675
$row213[0,7] = 0.033076663276222
$row213[0,8] = 0.047859454343135
$row213[0,9] = 0.030795744824741
$row213[0,10] = 0.024512227322229
$row213[0,11] = 0.01935006219525
$row213[0,12] = 0.061832623440327
$row213[0,13] = 0.02381220770047
$row213[0,14] = 0.024015797108296
$row213[0,15] = 0.026638685924874
$row213[0,16] = 0.025494557812284
$row213[0,17] = 0.023784432670898
$row213[0,18] = 0.023664239465129
$row213[0,19] = 0.024756608656363
$row213[0,20] = 0.028572937105183
$row213[0,21] = 0.034694583724618
$row213[0,22] = 0.028556239354123
$row213[0,23] = 0.034295528141972
$row213[0,24] = 0.029024798184609
$row213[0,25] = 0.02644593029774
$row213[0,26] = 0.025388330195042
$row213[0,27] = 0.07800007270701299
$row213[0,28] = 0.021637225940784
$row213[0,29] = 0.039404806859443
$row213[0,30] = 0.02735039896087
$row213[0,31] = 0.029494398053714
$row213[0,32] = 0.026301952845133
$row213[0,33] = 0.032198020095384
$row213[0,34] = 0.031908032236382
$row213[0,35] = 0.04536974368066
$row213[0,36] = 0.031656904707352
$row213[0,37] = 0.044210278895884
$row213[0,38] = 0.02998948509234
$row213[0,39] = 0.021670695112602
$row213[0,40] = 0.028443717166179
$row213[0,41] = 0.027692290788386
$row213[0,42] = 0.07255617000645
$row213[0,43] = 0.030869345280392
$row213[0,44] = 0.020813217719558
$row213[0,45] = 0.024346200482557
$row213[0,46] = 0.019584745764478
$row213[0,47] = 0.020163974340813
$row213[0,48] = 0.028372766754722
$row213[0,49] = 0.025394582567733
$row213[0,50] = 0.060531030972468
$row213[0,51] = 0.036699043919134
$row213[0,52] = 0.02915011816223
$row213[0,53] = 0.028459088967521
$row213[0,54] = 0.029858326074625
$row213[0,55] = 0.034006464720045
$ws.Range("B213:BE213").Value = $row213

$row214 = New-Object "object[,]" 1,56
$row214[0,0] = 0.080579780319622
$row214[0,1] = 0.072641303790099
$row214[0,2] = 0.078586818547809
$row214[0,3] = $null
$row214[0,4] = 0.062229474088845
$row214[0,5] = 0.056202581992138
$row214[0,6] = 0.058689849194675
$row214[0,7] = 0.051690302555801
$row214[0,8] = 0.04394065416683
$row214[0,9] = 0.047455380599396
$row214[0,10] = 0.053310034209468
$row214[0,11] = 0.05945427222032
$row214[0,12] = 0.058452370643448
$row214[0,13] = 0.032260699002383
$row214[0,14] = 0.07982375740641499
$row214[0,15] = 0.064873827425345
$row214[0,16] = 0.055453783662606
$row214[0,17] = 0.064952563125281
$row214[0,18] = 0.062882731268756
$row214[0,19] = 0.057366292785901
$row214[0,20] = 0.060812680115943
$row214[0,21] = 0.046559757775297
$row214[0,22] = 0.044071234694497
$row214[0,23] = 0.050857590883683
$row214[0,24] = 0.050744203432747
$row214[0,25] = 0.046763379216284
$row214[0,26] = 0.068998641051743
$row214[0,27] = 0.064804982579298
$row214[0,28] = 0.068826836475745
$row214[0,29] = 0.069501110291024
$row214[0,30] = 0.057736837850502
$row214[0,31] = 0.071657138370363
$row214[0,32] = 0.07935085395724301
$row214[0,33] = 0.055597846763709
$row214[0,34] = 0.046418791731054
$row214[0,35] = 0.05085369606832
$row214[0,36] = 0.049242888551309
$row214[0,37] = 0.047832180381964
$row214[0,38] = 0.061635805876434
$row214[0,39] = 0.06713622502963
$row214[0,40] = 0.048315008889494
$row214[0,41] = 0.056128119738376
$row214[0,42] = 0.07371135043476899
$row214[0,43] = 0.050442443690634
$row214[0,44] = 0.060930049349305
$row214[0,45] = 0.07499305889502
$row214[0,46] = 0.062206986008712
$row214[0,47] = 0.053154697984186
$row214[0,48] = 0.065097975908295
$row214[0,49] = 0.05061597276057
$row214[0,50] = 0.076243828428796
$row214[0,51] = 0.05578122938074
$row214[0,52] = 0.044387407960761
$row214[0,53] = 0.05158087884195
$row214[0,54] = 0.053890566359448
$row214[0,55] = 0.061931642116441
$ws.Range("B214:BE214").Value = $row214

$row215 = New-Object "object[,]" 1,56
$row215[0,0] = 0.048815118381829
$row215[0,1] = 0.068017980536516
$row215[0,2] = 0.06850193174315
$row215[0,3] = $null
$row215[0,4] = 0.042551845518156
$row215[0,5] = 0.051328740914253
$row215[0,6] = 0.064837132500166
$row215[0,7] = 0.059285520066162
$row215[0,8] = 0.054291935969209
$row215[0,9] = 0.057668463623795
$row215[0,10] = 0.06404348284788
$row215[0,11] = 0.06918160237580701
$row215[0,12] = 0.075242107444959
$row215[0,13] = 0.036235433800797
$row215[0,14] = 0.081633651734622
$row215[0,15] = 0.06968303413634
$row215[0,16] = 0.058902311621426
$row215[0,17] = 0.063443486726181
$row215[0,18] = 0.064493490799183
$row215[0,19] = 0.058937102728059
$row215[0,20] = 0.063928799924169
$row215[0,21] = 0.052948409224102
$row215[0,22] = 0.049074438141524
$row215[0,23] = 0.058564718073638
$row215[0,24] = 0.056986470106426
$row215[0,25] = 0.051778300380636
$row215[0,26] = 0.07397905540021101
$row215[0,27] = 0.11432176355704
$row215[0,28] = 0.07296188616244199
$row215[0,29] = 0.070221616936724
$row215[0,30] = 0.058919887611337
$row215[0,31] = 0.072545576557084
$row215[0,32] = 0.085778006332659
$row215[0,33] = 0.05903473569861
$row215[0,34] = 0.051675025131284
$row215[0,35] = 0.051586596734041
$row215[0,36] = 0.049435097156652
$row215[0,37] = 0.04960711616977
$row215[0,38] = 0.063902934658898
$row215[0,39] = 0.070104755470425
$row215[0,40] = 0.051069550701132
$row215[0,41] = 0.06246941211258
$row215[0,42] = 0.076996485402566
$row215[0,43] = 0.057280717387128
$row215[0,44] = 0.066347172488167
$row215[0,45] = 0.083223917181926
$row215[0,46] = 0.069345656071807
$row215[0,47] = 0.057823708089724
$row215[0,48] = 0.068623351820847
$row215[0,49] = 0.054708427564995
$row215[0,50] = 0.07870257154094
$row215[0,51] = 0.062368743095383
$row215[0,52] = 0.047343414271538
$row215[0,53] = 0.065409738139298
$row215[0,54] = 0.057785072085973
$row215[0,55] = 0.066905660951098
$ws.Range("B215:BE215").Value = $row215

$row216 = New-Object "object[,]" 1,56
$row216[0,0] = 0.063479621535476
$row216[0,1] = 0.06766739113063901
$row216[0,2] = 0.071937840180981
$row216[0,3] = $null
$row216[0,4] = 0.049271696353651
$row216[0,5] = 0.047334362713054
$row216[0,6] = 0.053671976569419
$row216[0,7] = 0.04888426187781
$row216[0,8] = 0.044491717981794
$row216[0,9] = 0.04596879139181
$row216[0,10] = 0.053657690131392
$row216[0,11] = 0.06012716839919
$row216[0,12] = 0.057490056554447
$row216[0,13] = 0.031427964666771
$row216[0,14] = 0.078021432912787
$row216[0,15] = 0.063695105971911
$row216[0,16] = 0.054378571729596
$row216[0,17] = 0.06419826437130501
$row216[0,18] = 0.06999636549549999
$row216[0,19] = 0.06336647314471
$row216[0,20] = 0.070167275946023
$row216[0,21] = 0.059857734935076
$row216[0,22] = 0.054395455917912
$row216[0,23] = 0.064622874064685
$row216[0,24] = 0.060304032555316
$row216[0,25] = 0.052877554672481
$row216[0,26] = 0.074499162991757
$row216[0,27] = 0.07014111590538299
$row216[0,28] = 0.0734831462237
$row216[0,29] = 0.068834398376065
$row216[0,30] = 0.055386915440452
$row216[0,31] = 0.069475729939623
$row216[0,32] = 0.085189941621674
$row216[0,33] = 0.060259536066108
$row216[0,34] = 0.049742596921268
$row216[0,35] = 0.05003144860697
$row216[0,36] = 0.04890632382611
$row216[0,37] = 0.051017789181242
$row216[0,38] = 0.060660911086842
$row216[0,39] = 0.064613831777953
$row216[0,40] = 0.046083858742904
$row216[0,41] = 0.053629985991348
$row216[0,42] = 0.078298298831025
$row216[0,43] = 0.046600887232026
$row216[0,44] = 0.057420277705546
$row216[0,45] = 0.074747252306528
$row216[0,46] = 0.060001815447056
$row216[0,47] = 0.052555806278147
$row216[0,48] = 0.0631846757184
$row216[0,49] = 0.049794937397444
$row216[0,50] = 0.071504295627267
$row216[0,51] = 0.054157219576257
$row216[0,52] = 0.046446818513254
$row216[0,53] = 0.064916380499527
$row216[0,54] = 0.056643918978411
$row216[0,55] = 0.069932171792205
$ws.Range("B216:BE216").Value = $row216

# --- Add new date label rows 217-221 (column A only) ---
$ws.Range("A216").Value = "02 09 2020"
$ws.Range("A217").Value = "03 09 2020"
$ws.Range("A218").Value = "04 09 2020"
$ws.Range("A219").Value = "05 09 2020"
$ws.Range("A220").Value = "06 09 2020"
$ws.Range("A221").Value = "07 09 2020"
